$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing A3 value: "RO.ACT.001.CRE" -> "RO.ACT.003HAB.SRA"
$ws.Range("A3").Value = "RO.ACT.003HAB.SRA"

# Add new value in A4: "RO.ACT.003HAB.SRL"
$ws.Range("A4").Value = "RO.ACT.003HAB.SRL"

# Update the active cell selection to B5
$ws.Range("B5").Select()
